$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1.34
$ws.Range("Q2").Value = 1.76
$ws.Range("S2").Value = 2.92
$ws.Range("X2").Value = 20
$ws.Range("F3").Value = 3.95
$ws.Range("G3").Value = 4.7
$ws.Range("H3").Value = 1.82
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 3.8
$ws.Range("K3").Value = 4.4
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 4.4
$ws.Range("O3").Value = 1.23
$ws.Range("P3").Value = 2.18
$ws.Range("Q3").Value = 1.68
$ws.Range("R3").Value = 1.46
$ws.Range("S3").Value = 2.74
$ws.Range("T3").Value = 1.66
$ws.Range("U3").Value = 2.24
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = 1.27
$ws.Range("X3").Value = 21
$ws.Range("Y3").Value = 12
$ws.Range("Z3").Value = 14.5
$ws.Range("AA3").Value = 24
$ws.Range("AB3").Value = 19
$ws.Range("AC3").Value = 9.800000000000001
$ws.Range("AD3").Value = 11
$ws.Range("AE3").Value = 20
$ws.Range("AG3").Value = 17.5
$ws.Range("AH3").Value = 17.5
$ws.Range("AO3").Value = 11.5
$ws.Range("N4").Value = 3
$ws.Range("P4").Value = 1.98
$ws.Range("Q4").Value = 1.84
$ws.Range("R4").Value = 1.31
$ws.Range("V4").Value = 1.18
$ws.Range("AH4").Value = 27
$ws.Range("S5").Value = 1.05
$ws.Range("G6").Value = 2.7
$ws.Range("I6").Value = 3.95
$ws.Range("P6").Value = 1.56
$ws.Range("R6").Value = 1.19
$ws.Range("S6").Value = 3.75
$ws.Range("W6").Value = 1.58
$ws.Range("AF8").Value = 70
$ws.Range("AH8").Value = 27
$ws.Range("M9").Value = 1.09
$ws.Range("N9").Value = 3.05
$ws.Range("Q9").Value = 2.22
$ws.Range("R9").Value = 1.26
$ws.Range("S9").Value = 4.3
$ws.Range("T9").Value = 2.16
$ws.Range("U9").Value = 1.74
$ws.Range("W9").Value = 2.38
$ws.Range("X9").Value = 11.5
$ws.Range("Y9").Value = 990
$ws.Range("Z9").Value = 60
$ws.Range("AA9").Value = 290
$ws.Range("AB9").Value = 6.8
$ws.Range("AC9").Value = 9
$ws.Range("AD9").Value = 990
$ws.Range("AE9").Value = 160
$ws.Range("AF9").Value = 9
$ws.Range("AG9").Value = 10.5
$ws.Range("AH9").Value = 990
$ws.Range("AI9").Value = 160
$ws.Range("AJ9").Value = 17
$ws.Range("AK9").Value = 980
$ws.Range("AL9").Value = 50
$ws.Range("AM9").Value = 240
$ws.Range("AN9").Value = 14
$ws.Range("F10").Value = 2.9
$ws.Range("G10").Value = 2.96
$ws.Range("L10").Value = 1.43
$ws.Range("W10").Value = 1.51
$ws.Range("G12").Value = 2.28
$ws.Range("N12").Value = 2.26
$ws.Range("O12").Value = 1.76
$ws.Range("W12").Value = 1.78
$ws.Range("AJ12").Value = 32
$ws.Range("F13").Value = 1.86
$ws.Range("I13").Value = 5.4
$ws.Range("V13").Value = 1.23
